$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E to Text format before writing, to preserve exact
# textual representation (e.g. "1.00", "42.973.10", percentages with padding)
# without Excel auto-converting them to numbers. ClearFormats afterwards so
# the cells end up with no explicit style, matching the original workbook.
$numRange = $ws.Range("D2:E51")
$numRange.NumberFormat = "@"

$ws.Range('D2').Value = '42.973.10'
$ws.Range('E2').Value = '  -0.89%  '
$ws.Range('D3').Value = '2.335.16'
$ws.Range('E3').Value = '  +1.25%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '306.53'
$ws.Range('E5').Value = '  -1.26%  '
$ws.Range('D6').Value = '100.73'
$ws.Range('E6').Value = '  -2.22%  '
$ws.Range('D7').Value = '0.512'
$ws.Range('E7').Value = '  -4.05%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').Value = '0.511'
$ws.Range('E9').Value = '  -3.95%  '
$ws.Range('D10').Value = '34.80'
$ws.Range('E10').Value = '  -2.36%  '
$ws.Range('D11').Value = '52.48'
$ws.Range('E11').Value = '  +0.84%  '
$ws.Range('D12').Value = '0.0797'
$ws.Range('E12').Value = '  -1.82%  '
$ws.Range('E13').Value = '  +0.71%  '
$ws.Range('D14').Value = '6.86'
$ws.Range('E14').Value = '  -2.18%  '
$ws.Range('D15').Value = '15.83'
$ws.Range('E15').Value = '  +5.73%  '
$ws.Range('D16').Value = '2.335.21'
$ws.Range('E16').Value = '  +2.40%  '
$ws.Range('D17').Value = '0.814'
$ws.Range('E17').Value = '  +0.78%  '
$ws.Range('D18').Value = '42.895.96'
$ws.Range('E18').Value = '  -0.83%  '
$ws.Range('D19').Value = '6.23'
$ws.Range('E19').Value = '  +0.95%  '
$ws.Range('D20').Value = '11.79'
$ws.Range('E20').Value = '  -3.76%  '
$ws.Range('D21').Value = '0.0₃0911'
$ws.Range('E21').Value = '  -2.47%  '
$ws.Range('D22').Value = '67.91'
$ws.Range('E22').Value = '  -0.28%  '
$ws.Range('D23').Value = '237.00'
$ws.Range('E23').Value = '  -1.87%  '
$ws.Range('D24').Value = '2.02'
$ws.Range('E24').Value = '  +0.71%  '
$ws.Range('E25').Value = '  -1.57%  '
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('D27').Value = '25.39'
$ws.Range('E27').Value = '  +1.62%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').Value = '2.32'
$ws.Range('E28').Value = '  +1.05%  '
$ws.Range('B29').Value = 'InjectiveProtocol'
$ws.Range('C29').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D29').Value = '35.27'
$ws.Range('E29').Value = '  -3.62%  '
$ws.Range('B30').Value = 'Cosmos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D30').Value = '9.40'
$ws.Range('E30').Value = '  -2.56%  '
$ws.Range('B31').Value = 'Monero'
$ws.Range('C31').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D31').Value = '163.62'
$ws.Range('E31').Value = '  -4.77%  '
$ws.Range('B32').Value = 'FirstDigitalUSD'
$ws.Range('C32').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D32').Value = '1.00'
$ws.Range('E32').Value = '  +0.03%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = '5.12'
$ws.Range('E33').Value = '  -2.96%  '
$ws.Range('B34').Value = 'Celestia'
$ws.Range('C34').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D34').Value = '17.52'
$ws.Range('E34').Value = '  -1.48%  '
$ws.Range('B35').Value = 'RenderToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D35').Value = '4.61'
$ws.Range('E35').Value = '  +6.58%  '
$ws.Range('E36').Value = '  -4.62%  '
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').Value = '0.0727'
$ws.Range('E37').Value = '  -1.84%  '
$ws.Range('B38').Value = 'ARBITRUM'
$ws.Range('C38').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D38').Value = '1.85'
$ws.Range('E38').Value = '  -1.36%  '
$ws.Range('B39').Value = 'LidoDAOToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D39').Value = '2.92'
$ws.Range('E39').Value = '  -4.77%  '
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').Value = '0.102'
$ws.Range('E40').Value = '  -3.48%  '
$ws.Range('B41').Value = 'Stellar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D41').Value = '0.113'
$ws.Range('E41').Value = '  -2.41%  '
$ws.Range('B42').Value = 'ApeXProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D42').Value = '2.58'
$ws.Range('E42').Value = '  +11.96%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '2.019.48'
$ws.Range('E43').Value = '  +2.43%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').Value = '0.0285'
$ws.Range('E44').Value = '  -2.61%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '18.77'
$ws.Range('E45').Value = '  -1.92%  '
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').Value = '10.17'
$ws.Range('E46').Value = '  +1.95%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').Value = '2.92'
$ws.Range('E47').Value = '  -2.36%  '
$ws.Range('B48').Value = 'MultiversX'
$ws.Range('C48').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D48').Value = '56.12'
$ws.Range('E48').Value = '  +1.22%  '
$ws.Range('B49').Value = 'HuobiToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D49').Value = '2.90'
$ws.Range('E49').Value = '  -1.28%  '
$ws.Range('B50').Value = 'RocketPoolETH'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D50').Value = '2.559.15'
$ws.Range('E50').Value = '  +1.09%  '
$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D51').Value = '4.67'
$ws.Range('E51').Value = '  +1.55%  '

$numRange.ClearFormats()

